# Update the account-example worksheet:
#  - the "matKhau" (password) column is no longer returned by the
#    get-user-list endpoint, so drop it from the sample/import sheet
#  - row 2/3 data shifts up one column (C->B, D->C)
#  - the old B2 hyperlink (mailto:) lived on the password column and
#    goes away with it
#  - selection moves to G10 (matches the new "import users from file"
#    sample focus)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the hyperlink that lived on B2 (password column) before the
# column shift carries it along to the wrong cell.
$ws.Hyperlinks.Delete()

# The hyperlink cell format (s="1", "Hyperlink" cell style/font) must
# go together with the hyperlink itself.
$ws.Range("B2").ClearFormats()
$wb.Styles.Item("Hyperlink").Delete()

# Remove the whole "matKhau" column; loai/trangThai shift left.
$ws.Columns("B").Delete()

# Match the new column widths used by the updated sample file
# (closest values the host's pixel-quantized width model can hit).
$ws.Columns("A").ColumnWidth = 19.25
$ws.Columns("B").ColumnWidth = 14.75

# Selection as left by the last save.
$ws.Range("G10").Select() | Out-Null
